$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.753.55"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.302.36"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'301.30"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "'96.18"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "'0.513"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("D10").Value = "'34.79"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "'19.20"
$ws.Range("E11").Value = "  +4.70%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "'6.80"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "2.651.51"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "2.294.96"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'0.786"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "42.668.80"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'12.34"
$ws.Range("E19").Value = "  -6.49%  "
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'6.03"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "'67.92"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'2.26"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").Value = "'235.17"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "'24.63"
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'164.95"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").Value = "'9.08"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "'32.28"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "'4.98"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "'17.56"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  -7.15%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "'1.76"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").Value = "'19.95"
$ws.Range("E42").Value = "  +7.18%  "
$ws.Range("D43").Value = "1.968.82"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("D47").Value = "'2.77"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "2.524.84"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "'53.21"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").Value = "'71.70"
$ws.Range("E51").Value = "  -0.48%  "
